$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New vocabulary rows appended below the existing table (rows 57-60), plus
# filler rows 61-68 which keep the "=" helper formulas in sync with the rest
# of the sheet (mirrors the existing pattern used for previously-empty rows).
# ---------------------------------------------------------------------------

# Row 57: percent / Prozent
$ws.Range("A57").Value = "percent"
$ws.Range("B57").Value = "Prozent"
$ws.Range("C57").Formula = '=A57&"="&A57'
$ws.Range("D57").Formula = '=A57&"="&B57'
$ws.Range("E57").Formula = '=SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(D57,"ß","ss"),"ä","ae"),"ö","oe"),"ü","ue")'
$ws.Range("F57").Formula = '=(SUBSTITUTE(C57," ","\u0020"))'
$ws.Range("G57").Formula = '=(SUBSTITUTE(D57," ","\u0020"))'

# Row 58: Processing / Verarbeitung der Daten (German cell uses the small Arial style)
$ws.Range("A58").Value = "Processing"
$ws.Range("B58").Value = "Verarbeitung der Daten"
$ws.Range("B58").Font.Name = "Arial"
$ws.Range("B58").Font.Size = 10
$ws.Range("B58").Font.Color = 0
$ws.Range("C58").Formula = '=A58&"="&A58'
$ws.Range("D58").Formula = '=A58&"="&B58'
$ws.Range("E58").Formula = '=SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(D58,"ß","ss"),"ä","ae"),"ö","oe"),"ü","ue")'
$ws.Range("F58").Formula = '=(SUBSTITUTE(C58," ","\u0020"))'
$ws.Range("G58").Formula = '=(SUBSTITUTE(D58," ","\u0020"))'

# Row 59: Searching for music video files / Suche nach Musikvideos
$ws.Range("A59").Value = "Searching for music video files"
$ws.Range("B59").Value = "Suche nach Musikvideos"

# Row 60: note the German text is entered before the English one, matching
# the shared-string order captured by the original authoring session.
$ws.Range("B60").Value = "Es müssen zuerst Musikvideos hinzugefügt werden um diese Funktion nützen zu können"
$ws.Range("A60").Value = "You first need to add music videos to use this feature"

# Rows 59-68 share one formula block per column (new shared-formula group).
$ws.Range("C59:C68").Formula = '=A59&"="&A59'
$ws.Range("D59:D68").Formula = '=A59&"="&B59'
$ws.Range("E59:E68").Formula = '=SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(D59,"ß","ss"),"ä","ae"),"ö","oe"),"ü","ue")'
$ws.Range("F59:F68").Formula = '=(SUBSTITUTE(C59," ","\u0020"))'
$ws.Range("G59:G68").Formula = '=(SUBSTITUTE(D59," ","\u0020"))'

# ---------------------------------------------------------------------------
# Grow the table (ListObject) so the new rows belong to "Tabelle1".
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:B67"))

# ---------------------------------------------------------------------------
# Conditional formatting on column D now covers the extended range.
# ---------------------------------------------------------------------------
$ws.Range("D2:D68").FormatConditions.Delete()
$cf = $ws.Range("D2:D68").FormatConditions.Add(2, 0, 'SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(SUBSTITUTE(#REF!, "ü","ue"),"ö","oe"),"ä","ae"),"Ü","Ue"),"Ö","Oe"),"Ä","Ae"),"ß","ss")')

# ---------------------------------------------------------------------------
# Restore the view: scroll position + active cell selection.
# ---------------------------------------------------------------------------
$ws.Range("C73").Select()
